# Weekly fruit/vegetable price update: insert a new observation row for
# "Vega Modelo de Temuco" / Zapallo (Camote, 1a nueva(o), origin Perú)
# dated 2021-11-16 (serial 44516) at row 391, pushing the existing rows
# 391-416 down to 392-417.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 391 (shifts rows 391..416 -> 392..417)
$ws.Rows.Item(391).Insert()

# Populate the new row with the new data point
$ws.Cells.Item(391, 1).Value  = 10
$ws.Cells.Item(391, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(391, 3).Value  = "La Araucanía"
$ws.Cells.Item(391, 4).Value  = 44516
$ws.Cells.Item(391, 5).Value  = 9
$ws.Cells.Item(391, 6).Value  = 100112045
$ws.Cells.Item(391, 7).Value  = "Zapallo"
$ws.Cells.Item(391, 8).Value  = "Camote"
$ws.Cells.Item(391, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(391, 10).Value = 650
$ws.Cells.Item(391, 11).Value = 800
$ws.Cells.Item(391, 12).Value = 800
$ws.Cells.Item(391, 13).Value = 800
$ws.Cells.Item(391, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(391, 15).Value = "Perú"
$ws.Cells.Item(391, 16).Value = 800
$ws.Cells.Item(391, 17).Value = 1
$ws.Cells.Item(391, 18).Value = "Hortaliza"
